$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.821.64'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '2.470.82'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.09'
$ws.Range('C5').Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4122) | Out-Null
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.26'
$ws.Range('C6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4122) | Out-Null
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('C8').Copy() | Out-Null
$ws.Range('D8').PasteSpecial(-4122) | Out-Null
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('E9').Value = '  +5.64%  '
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.85'
$ws.Range('C12').Copy() | Out-Null
$ws.Range('D12').PasteSpecial(-4122) | Out-Null
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').Value = '68.761.33'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.62'
$ws.Range('C15').Copy() | Out-Null
$ws.Range('D15').PasteSpecial(-4122) | Out-Null
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '10.65'
$ws.Range('C16').Copy() | Out-Null
$ws.Range('D16').PasteSpecial(-4122) | Out-Null
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '338.45'
$ws.Range('C17').Copy() | Out-Null
$ws.Range('D17').PasteSpecial(-4122) | Out-Null
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.94'
$ws.Range('C18').Copy() | Out-Null
$ws.Range('D18').PasteSpecial(-4122) | Out-Null
$ws.Range('E18').Value = '  -2.88%  '
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('C20').Copy() | Out-Null
$ws.Range('D20').PasteSpecial(-4122) | Out-Null
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.88'
$ws.Range('C21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4122) | Out-Null
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.97'
$ws.Range('C22').Copy() | Out-Null
$ws.Range('D22').PasteSpecial(-4122) | Out-Null
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.67'
$ws.Range('C23').Copy() | Out-Null
$ws.Range('D23').PasteSpecial(-4122) | Out-Null
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.26'
$ws.Range('C24').Copy() | Out-Null
$ws.Range('D24').PasteSpecial(-4122) | Out-Null
$ws.Range('E24').Value = '  +1.97%  '
$ws.Range('D25').Value = '0.0₃0825'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.24'
$ws.Range('C26').Copy() | Out-Null
$ws.Range('D26').PasteSpecial(-4122) | Out-Null
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '430.41'
$ws.Range('C28').Copy() | Out-Null
$ws.Range('D28').PasteSpecial(-4122) | Out-Null
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('E30').Value = '  -2.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.81'
$ws.Range('C31').Copy() | Out-Null
$ws.Range('D31').PasteSpecial(-4122) | Out-Null
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.00'
$ws.Range('C32').Copy() | Out-Null
$ws.Range('D32').PasteSpecial(-4122) | Out-Null
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('C34').Copy() | Out-Null
$ws.Range('D34').PasteSpecial(-4122) | Out-Null
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.88'
$ws.Range('C35').Copy() | Out-Null
$ws.Range('D35').PasteSpecial(-4122) | Out-Null
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.46'
$ws.Range('C36').Copy() | Out-Null
$ws.Range('D36').PasteSpecial(-4122) | Out-Null
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.299'
$ws.Range('C37').Copy() | Out-Null
$ws.Range('D37').PasteSpecial(-4122) | Out-Null
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.08'
$ws.Range('C39').Copy() | Out-Null
$ws.Range('D39').PasteSpecial(-4122) | Out-Null
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.07'
$ws.Range('C40').Copy() | Out-Null
$ws.Range('D40').PasteSpecial(-4122) | Out-Null
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.39'
$ws.Range('C41').Copy() | Out-Null
$ws.Range('D41').PasteSpecial(-4122) | Out-Null
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '130.67'
$ws.Range('C42').Copy() | Out-Null
$ws.Range('D42').PasteSpecial(-4122) | Out-Null
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0720'
$ws.Range('C43').Copy() | Out-Null
$ws.Range('D43').PasteSpecial(-4122) | Out-Null
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.487'
$ws.Range('C44').Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4122) | Out-Null
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.565'
$ws.Range('C45').Copy() | Out-Null
$ws.Range('D45').PasteSpecial(-4122) | Out-Null
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0920'
$ws.Range('C46').Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4122) | Out-Null
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.01'
$ws.Range('C49').Copy() | Out-Null
$ws.Range('D49').PasteSpecial(-4122) | Out-Null
$ws.Range('E49').Value = '  -6.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.91'
$ws.Range('C50').Copy() | Out-Null
$ws.Range('D50').PasteSpecial(-4122) | Out-Null
$ws.Range('E50').Value = '  -3.25%  '
$ws.Range('D51').Value = '0.0₆0207'
$ws.Range('E51').Value = '  -8.75%  '

$excel.CutCopyMode = 0
